# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet to reflect the latest scraped values, mirroring the change
# produced by the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new Price (column D) / new Volume(1h) (column E) text.
# $null means that column is unchanged for that row.
$updates = @(
    @{Row=2; D='25.800.28'; E='  +0.23%  '},
    @{Row=3; D='1.737.83'; E='  -0.71%  '},
    @{Row=4; D='1.001'; E='  -0.13%  '},
    @{Row=5; D='226.10'; E='  -4.14%  '},
    @{Row=6; D=$null; E='  -0.09%  '},
    @{Row=7; D='0.5150'; E='  +2.00%  '},
    @{Row=8; D='0.2720'; E='  +3.69%  '},
    @{Row=9; D='38.65'; E='  -5.92%  '},
    @{Row=10; D='0.06081'; E='  -2.29%  '},
    @{Row=11; D='1.740.63'; E='  -0.52%  '},
    @{Row=12; D='0.06994'; E='  +1.13%  '},
    @{Row=13; D='15.22'; E='  -1.32%  '},
    @{Row=14; D='0.6303'; E='  +5.86%  '},
    @{Row=15; D='4.485'; E='  +0.39%  '},
    @{Row=16; D='76.22'; E='  -0.48%  '},
    @{Row=17; D='1.001'; E='  -0.11%  '},
    @{Row=18; D='0.9999'; E='  -0.18%  '},
    @{Row=19; D='25.831.17'; E='  +0.30%  '},
    @{Row=20; D='11.41'; E='  -1.47%  '},
    @{Row=21; D='0.000006567'; E='  -3.23%  '},
    @{Row=22; D='1.959.80'; E='  -0.67%  '},
    @{Row=23; D=$null; E='  -0.67%  '},
    @{Row=24; D='8.413'; E='  +2.45%  '},
    @{Row=25; D='5.088'; E='  -1.48%  '},
    @{Row=26; D='136.01'; E='  -0.51%  '},
    @{Row=27; D=$null; E='  +4.08%  '},
    @{Row=28; D='1.812'; E='  +0.71%  '},
    @{Row=29; D='14.95'; E='  +0.32%  '},
    @{Row=30; D='102.39'; E='  +0.39%  '},
    @{Row=31; D='0.08305'; E='  +1.76%  '},
    @{Row=32; D='3.604'; E='  -1.44%  '},
    @{Row=33; D='3.357'; E='  -1.60%  '},
    @{Row=34; D=$null; E='  -1.56%  '},
    @{Row=35; D='2.609'; E='  -1.84%  '},
    @{Row=36; D='0.9708'; E='  -1.75%  '},
    @{Row=37; D='0.5936'; E='  -1.55%  '},
    @{Row=38; D='2.681'; E='  -0.10%  '},
    @{Row=40; D='1.939'; E='  +0.57%  '},
    @{Row=41; D='0.9988'; E='  -0.24%  '},
    @{Row=42; D='101.79'; E='  -0.74%  '},
    @{Row=43; D='0.3782'; E='  +0.24%  '},
    @{Row=44; D='0.7227'; E='  -1.75%  '},
    @{Row=45; D='4.856'; E='  -1.37%  '},
    @{Row=46; D='0.05486'; E='  +0.24%  '},
    @{Row=47; D='6.221'; E='  +5.55%  '},
    @{Row=48; D='0.1097'; E='  +0.38%  '},
    @{Row=49; D='29.70'; E='  +0.14%  '},
    @{Row=50; D='51.78'; E='  -0.39%  '},
    @{Row=51; D='1.001'; E='  +0.09%  '}
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D
        # Force text so values like "1.001" aren't reinterpreted as numbers,
        # then restore the default style so no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($row, 5)   # column E
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}

$wb.Save()
